# Fruta / hortaliza, semanal
# Insert a new weekly record at row 122 (pushing the existing rows 122:151
# down to 123:152) and populate it with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 122:151 down to 123:152, preserving row 122's original layout.
$ws.Rows.Item(122).Insert()

$newRow = 122

$ws.Cells.Item($newRow, 1).Value2  = 10
$ws.Cells.Item($newRow, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item($newRow, 3).Value2  = "La Araucanía"
$ws.Cells.Item($newRow, 4).Value2  = 45275
$ws.Cells.Item($newRow, 5).Value2  = 9
$ws.Cells.Item($newRow, 6).Value2  = "Fruta"
$ws.Cells.Item($newRow, 7).Value2  = 100108
$ws.Cells.Item($newRow, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item($newRow, 9).Value2  = 100108004
$ws.Cells.Item($newRow, 10).Value2 = "Papaya"
$ws.Cells.Item($newRow, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item($newRow, 12).Value2 = "Primera"
$ws.Cells.Item($newRow, 13).Value2 = 40
$ws.Cells.Item($newRow, 14).Value2 = 24000
$ws.Cells.Item($newRow, 15).Value2 = 24000
$ws.Cells.Item($newRow, 16).Value2 = 24000
$ws.Cells.Item($newRow, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item($newRow, 18).Value2 = "Provincia del Elquí"
$ws.Cells.Item($newRow, 19).Value2 = 2400
$ws.Cells.Item($newRow, 20).Value2 = 10
